$wb = $excel.ActiveWorkbook

# --- Sheet "info" updates ---
$info = $wb.Worksheets.Item("info")
# A1 holds a numeric-looking id that must stay text, like the original "0001"
$info.Range("A1").NumberFormat = "@"
$info.Range("A1").Value = "99393"
$info.Range("B1").Value = "dk"
$info.Range("C1").Value = "dk"

# --- Sheet "items" updates ---
$items = $wb.Worksheets.Item("items")
$items.Range("A1").Value = "사이다"
$items.Range("C1").Value = 1000
$items.Range("E1").Value = 1000

$items.Range("A2").Value = "소주"
$items.Range("B2").Value = "개"
$items.Range("C2").Value = 4000
$items.Range("D2").Value = 1
$items.Range("E2").Value = 4000

$items.Range("A3").Value = "편육"
$items.Range("B3").Value = "개"
$items.Range("C3").Value = 12000
$items.Range("D3").Value = 1
$items.Range("E3").Value = 12000
